# cs-en-us-112pct.xlsx weekly refresh: "New crime data collected"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (rich-text shared strings): volume/number + reporting week
# ---------------------------------------------------------------------------
# "Volume 32   Number  32"  ->  "Volume 32   Number  34"
$ws.Range("A8").Value = "Volume 32   Number  34"

# "Report Covering the Week  8/4/2025  Through  8/10/2025"
#   -> "Report Covering the Week  8/18/2025  Through  8/24/2025"
$ws.Range("C9").Value = "Report Covering the Week  8/18/2025  Through  8/24/2025"

# ---------------------------------------------------------------------------
# Cells that flip between the "no data" text placeholders ("0" / "***.*")
# and real numbers need both their style AND value fixed up. Copying from a
# cell that already carries the desired style+content reproduces the exact
# style index instead of minting a new one.
# ---------------------------------------------------------------------------
$ws.Range("F15").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2

$ws.Range("F15").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2

$ws.Range("K15").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

$ws.Range("F15").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2

$ws.Range("C15").Copy($ws.Range("D27"))   # -> text "0"
$ws.Range("E15").Copy($ws.Range("E27"))   # -> text "***.*"

$ws.Range("C15").Copy($ws.Range("C28"))   # -> text "0"

# ---------------------------------------------------------------------------
# Row 15 - Murder
# ---------------------------------------------------------------------------
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = -25

# ---------------------------------------------------------------------------
# Row 16 - Rape
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -77.777777777777
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = -12.5
$ws.Range("L16").Value = 13.953488372093
$ws.Range("M16").Value = -28.985507246376
$ws.Range("N16").Value = -88.683602771362

# ---------------------------------------------------------------------------
# Row 17 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 109
$ws.Range("J17").Value = 69
$ws.Range("K17").Value = 57.971014492753
$ws.Range("L17").Value = 65.151515151515
$ws.Range("M17").Value = 194.594594594595
$ws.Range("N17").Value = 36.25

# ---------------------------------------------------------------------------
# Row 18 - Fel. Assault  (C18 style/value already fixed above)
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = 39.705882352941
$ws.Range("L18").Value = 23.376623376623
$ws.Range("M18").Value = 11.764705882352
$ws.Range("N18").Value = -89.030023094688

# ---------------------------------------------------------------------------
# Row 19 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 236
$ws.Range("J19").Value = 297
$ws.Range("K19").Value = -20.538720538720
$ws.Range("L19").Value = -27.160493827160
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -63.354037267080

# ---------------------------------------------------------------------------
# Row 20 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 103
$ws.Range("J20").Value = 112
$ws.Range("K20").Value = -8.035714285714
$ws.Range("L20").Value = 9.574468085106
$ws.Range("M20").Value = 58.461538461538
$ws.Range("N20").Value = -95.326678765880

# ---------------------------------------------------------------------------
# Row 21 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -30.434782608695
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 77
$ws.Range("H21").Value = -18.181818181818
$ws.Range("I21").Value = 599
$ws.Range("J21").Value = 607
$ws.Range("K21").Value = -1.317957166392
$ws.Range("L21").Value = -1.317957166392
$ws.Range("M21").Value = 21.010101010101
$ws.Range("N21").Value = -85.865974516281

# ---------------------------------------------------------------------------
# Row 22 - TOTAL  (D22/E22 style/value already fixed above)
# ---------------------------------------------------------------------------
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = -24
$ws.Range("L22").Value = -9.523809523809
$ws.Range("M22").Value = 46.153846153846

# ---------------------------------------------------------------------------
# Row 24 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = -27.777777777777
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 174
$ws.Range("H24").Value = -18.390804597701
$ws.Range("I24").Value = 1282
$ws.Range("J24").Value = 1107
$ws.Range("K24").Value = 15.808491418247
$ws.Range("L24").Value = 29.888551165146
$ws.Range("M24").Value = 99.068322981366

# ---------------------------------------------------------------------------
# Row 25 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 38
$ws.Range("E25").Value = -21.052631578947
$ws.Range("F25").Value = 108
$ws.Range("G25").Value = 121
$ws.Range("H25").Value = -10.743801652892
$ws.Range("I25").Value = 999
$ws.Range("J25").Value = 797
$ws.Range("K25").Value = 25.345043914680
$ws.Range("L25").Value = 44.573082489146

# ---------------------------------------------------------------------------
# Row 26 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 188
$ws.Range("J26").Value = 183
$ws.Range("K26").Value = 2.732240437158
$ws.Range("L26").Value = 22.077922077922
$ws.Range("M26").Value = 33.333333333333

# ---------------------------------------------------------------------------
# Row 27 - Misd. Assault  (C27/D27/E27 style/value already fixed above)
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 13
$ws.Range("K27").Value = 8.333333333333
$ws.Range("L27").Value = 62.5

# ---------------------------------------------------------------------------
# Row 28 - UCR Rape*  (C28 style/value already fixed above)
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 19.047619047619
$ws.Range("L28").Value = 19.047619047619

# ---------------------------------------------------------------------------
# Row 31 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("I31").Value = 10
$ws.Range("K31").Value = 100
$ws.Range("L31").Value = -16.666666666666
